$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1) Row 1546 was previously the last row of the response table, so it
#    carried the special "bottom border" style (19/20/21). New
#    responses are being appended below it, so it becomes a normal
#    alternating-color row. It sits at an even offset in the table body
#    (header is row 1), so it should now look like row 1544, the
#    preceding even-styled row. Only the formatting changes -- the
#    timestamp and every shared-string reference stay exactly as they
#    were.
# ---------------------------------------------------------------------
$ws.Range('A1544:H1544').Copy()
$ws.Range('A1546:H1546').PasteSpecial(-4122)
$excel.CutCopyMode = 0

# ---------------------------------------------------------------------
# 2) Append the 21 new form responses (rows 1547-1567). Each new row's
#    look-and-feel is cloned from an existing row that already has the
#    right alternating style (and, where needed, the special
#    "adolorido" highlight on column G) via Copy, then the actual
#    values for that response are written on top.
# ---------------------------------------------------------------------

$ws.Range('A1545:H1545').Copy($ws.Range('A1547:H1547'))
$ws.Cells.Item(1547,1).Value2 = 45880.351948807875
$ws.Cells.Item(1547,2).Value2 = 'Lunes'
$ws.Cells.Item(1547,3).Value2 = 'Alejandro Zendejas'
$ws.Cells.Item(1547,4).Value2 = 'Normal'
$ws.Cells.Item(1547,5).Value2 = 'Mejor que normal'
$ws.Cells.Item(1547,6).Value2 = 'Más de 8'
$ws.Cells.Item(1547,7).Value2 = 'Normal'
$ws.Cells.Item(1547,8).Value2 = 'Nada'

$ws.Range('A1544:H1544').Copy($ws.Range('A1548:H1548'))
$ws.Cells.Item(1548,1).Value2 = 45880.35497219907
$ws.Cells.Item(1548,2).Value2 = 'Lunes'
$ws.Cells.Item(1548,3).Value2 = 'Néstor Araujo'
$ws.Cells.Item(1548,4).Value2 = 'Normal'
$ws.Cells.Item(1548,5).Value2 = 'Normal'
$ws.Cells.Item(1548,6).Value2 = '6 a 8'
$ws.Cells.Item(1548,7).Value2 = 'Normal'
$ws.Cells.Item(1548,8).Value2 = 'Nada'

$ws.Range('A1545:H1545').Copy($ws.Range('A1549:H1549'))
$ws.Cells.Item(1549,1).Value2 = 45880.35605615741
$ws.Cells.Item(1549,2).Value2 = 'Lunes'
$ws.Cells.Item(1549,3).Value2 = 'José Raúl Zúñiga'
$ws.Cells.Item(1549,4).Value2 = 'Normal'
$ws.Cells.Item(1549,5).Value2 = 'Peor que normal'
$ws.Cells.Item(1549,6).Value2 = '6 a 8'
$ws.Cells.Item(1549,7).Value2 = 'Normal'
$ws.Cells.Item(1549,8).Value2 = 'Nada'

$ws.Range('A1544:H1544').Copy($ws.Range('A1550:H1550'))
$ws.Cells.Item(1550,1).Value2 = 45880.3568934375
$ws.Cells.Item(1550,2).Value2 = 'Lunes'
$ws.Cells.Item(1550,3).Value2 = 'Luis Ángel Malagón'
$ws.Cells.Item(1550,4).Value2 = 'Normal'
$ws.Cells.Item(1550,5).Value2 = 'Normal'
$ws.Cells.Item(1550,6).Value2 = '6 a 8'
$ws.Cells.Item(1550,7).Value2 = 'Normal'
$ws.Cells.Item(1550,8).Value2 = 'Nada'

$ws.Range('A1545:H1545').Copy($ws.Range('A1551:H1551'))
$ws.Cells.Item(1551,1).Value2 = 45880.36914304398
$ws.Cells.Item(1551,2).Value2 = 'Lunes'
$ws.Cells.Item(1551,3).Value2 = 'Israel Reyes'
$ws.Cells.Item(1551,4).Value2 = 'Normal'
$ws.Cells.Item(1551,5).Value2 = 'Peor que normal'
$ws.Cells.Item(1551,6).Value2 = 'Menos de 6'
$ws.Cells.Item(1551,7).Value2 = 'Normal'
$ws.Cells.Item(1551,8).Value2 = 'Nada'

$ws.Range('A1544:H1544').Copy($ws.Range('A1552:H1552'))
$ws.Cells.Item(1552,1).Value2 = 45880.37883199074
$ws.Cells.Item(1552,2).Value2 = 'Lunes'
$ws.Cells.Item(1552,3).Value2 = 'Jonathan Dos Santos'
$ws.Cells.Item(1552,4).Value2 = 'Normal'
$ws.Cells.Item(1552,5).Value2 = 'Normal'
$ws.Cells.Item(1552,6).Value2 = '6 a 8'
$ws.Cells.Item(1552,7).Value2 = 'Normal'
$ws.Cells.Item(1552,8).Value2 = 'Nada'

$ws.Range('A1545:H1545').Copy($ws.Range('A1553:H1553'))
$ws.Cells.Item(1553,1).Value2 = 45880.38362291666
$ws.Cells.Item(1553,2).Value2 = 'Lunes'
$ws.Cells.Item(1553,3).Value2 = 'Alexis Gutiérrez'
$ws.Cells.Item(1553,4).Value2 = 'Normal'
$ws.Cells.Item(1553,5).Value2 = 'Normal'
$ws.Cells.Item(1553,6).Value2 = '6 a 8'
$ws.Cells.Item(1553,7).Value2 = 'Normal'
$ws.Cells.Item(1553,8).Value2 = 'Nada'

$ws.Range('A1544:H1544').Copy($ws.Range('A1554:H1554'))
$ws.Cells.Item(1554,1).Value2 = 45880.384631030094
$ws.Cells.Item(1554,2).Value2 = 'Lunes'
$ws.Cells.Item(1554,3).Value2 = 'Erick Sánchez'
$ws.Cells.Item(1554,4).Value2 = 'Normal'
$ws.Cells.Item(1554,5).Value2 = 'Normal'
$ws.Cells.Item(1554,6).Value2 = '6 a 8'
$ws.Cells.Item(1554,7).Value2 = 'Normal'
$ws.Cells.Item(1554,8).Value2 = 'Nada'

$ws.Range('A1545:H1545').Copy($ws.Range('A1555:H1555'))
$ws.Cells.Item(1555,1).Value2 = 45880.38479652778
$ws.Cells.Item(1555,2).Value2 = 'Lunes'
$ws.Cells.Item(1555,3).Value2 = 'Brian Rodríguez'
$ws.Cells.Item(1555,4).Value2 = 'Normal'
$ws.Cells.Item(1555,5).Value2 = 'Normal'
$ws.Cells.Item(1555,6).Value2 = '6 a 8'
$ws.Cells.Item(1555,7).Value2 = 'Normal'
$ws.Cells.Item(1555,8).Value2 = 'Nada'

$ws.Range('A1544:H1544').Copy($ws.Range('A1556:H1556'))
$ws.Cells.Item(1556,1).Value2 = 45880.38497909722
$ws.Cells.Item(1556,2).Value2 = 'Lunes'
$ws.Cells.Item(1556,3).Value2 = 'Rodrigo Aguirre'
$ws.Cells.Item(1556,4).Value2 = 'Normal'
$ws.Cells.Item(1556,5).Value2 = 'Normal'
$ws.Cells.Item(1556,6).Value2 = '6 a 8'
$ws.Cells.Item(1556,7).Value2 = 'Normal'
$ws.Cells.Item(1556,8).Value2 = 'Nada'

$ws.Range('A139:H139').Copy($ws.Range('A1557:H1557'))
$ws.Cells.Item(1557,1).Value2 = 45880.38532104167
$ws.Cells.Item(1557,2).Value2 = 'Lunes'
$ws.Cells.Item(1557,3).Value2 = 'Álvaro Fidalgo'
$ws.Cells.Item(1557,4).Value2 = 'Normal'
$ws.Cells.Item(1557,5).Value2 = 'Normal'
$ws.Cells.Item(1557,6).Value2 = '6 a 8'
$ws.Cells.Item(1557,7).Value2 = 'Adolorido de una zona'
$ws.Cells.Item(1557,8).Value2 = '9 Isquiotibial izquierdo'

$ws.Range('A1544:H1544').Copy($ws.Range('A1558:H1558'))
$ws.Cells.Item(1558,1).Value2 = 45880.387728368056
$ws.Cells.Item(1558,2).Value2 = 'Lunes'
$ws.Cells.Item(1558,3).Value2 = 'Santiago Naveda'
$ws.Cells.Item(1558,4).Value2 = 'Normal'
$ws.Cells.Item(1558,5).Value2 = 'Normal'
$ws.Cells.Item(1558,6).Value2 = '6 a 8'
$ws.Cells.Item(1558,7).Value2 = 'Normal'
$ws.Cells.Item(1558,8).Value2 = 'Nada'

$ws.Range('A139:H139').Copy($ws.Range('A1559:H1559'))
$ws.Cells.Item(1559,1).Value2 = 45880.38928990741
$ws.Cells.Item(1559,2).Value2 = 'Lunes'
$ws.Cells.Item(1559,3).Value2 = 'Henry Martín'
$ws.Cells.Item(1559,4).Value2 = 'Normal'
$ws.Cells.Item(1559,5).Value2 = 'Normal'
$ws.Cells.Item(1559,6).Value2 = '6 a 8'
$ws.Cells.Item(1559,7).Value2 = 'Adolorido de una zona'
$ws.Cells.Item(1559,8).Value2 = '9 Isquiotibial izquierdo'

$ws.Range('A1544:H1544').Copy($ws.Range('A1560:H1560'))
$ws.Cells.Item(1560,1).Value2 = 45880.389651053236
$ws.Cells.Item(1560,2).Value2 = 'Lunes'
$ws.Cells.Item(1560,3).Value2 = 'Sebastián Cáceres'
$ws.Cells.Item(1560,4).Value2 = 'Normal'
$ws.Cells.Item(1560,5).Value2 = 'Normal'
$ws.Cells.Item(1560,6).Value2 = '6 a 8'
$ws.Cells.Item(1560,7).Value2 = 'Normal'
$ws.Cells.Item(1560,8).Value2 = 'Nada'

$ws.Range('A1545:H1545').Copy($ws.Range('A1561:H1561'))
$ws.Cells.Item(1561,1).Value2 = 45880.39028413194
$ws.Cells.Item(1561,2).Value2 = 'Lunes'
$ws.Cells.Item(1561,3).Value2 = 'Erick Sánchez'
$ws.Cells.Item(1561,4).Value2 = 'Normal'
$ws.Cells.Item(1561,5).Value2 = 'Normal'
$ws.Cells.Item(1561,6).Value2 = '6 a 8'
$ws.Cells.Item(1561,7).Value2 = 'Normal'
$ws.Cells.Item(1561,8).Value2 = 'Nada'

$ws.Range('A1460:H1460').Copy($ws.Range('A1562:H1562'))
$ws.Cells.Item(1562,1).Value2 = 45880.39076241898
$ws.Cells.Item(1562,2).Value2 = 'Lunes'
$ws.Cells.Item(1562,3).Value2 = 'Isaías Violante'
$ws.Cells.Item(1562,4).Value2 = 'Normal'
$ws.Cells.Item(1562,5).Value2 = 'Normal'
$ws.Cells.Item(1562,6).Value2 = '6 a 8'
$ws.Cells.Item(1562,7).Value2 = 'Adolorido de una zona'
$ws.Cells.Item(1562,8).Value2 = '19 hombro derecho'

$ws.Range('A139:H139').Copy($ws.Range('A1563:H1563'))
$ws.Cells.Item(1563,1).Value2 = 45880.39178005787
$ws.Cells.Item(1563,2).Value2 = 'Lunes'
$ws.Cells.Item(1563,3).Value2 = 'Cristian Borja'
$ws.Cells.Item(1563,4).Value2 = 'Normal'
$ws.Cells.Item(1563,5).Value2 = 'Normal'
$ws.Cells.Item(1563,6).Value2 = '6 a 8'
$ws.Cells.Item(1563,7).Value2 = 'Adolorido de una zona'
$ws.Cells.Item(1563,8).Value2 = '16 espalda'

$ws.Range('A1460:H1460').Copy($ws.Range('A1564:H1564'))
$ws.Cells.Item(1564,1).Value2 = 45880.39196759259
$ws.Cells.Item(1564,2).Value2 = 'Lunes'
$ws.Cells.Item(1564,3).Value2 = 'Víctor Dávila'
$ws.Cells.Item(1564,4).Value2 = 'Normal'
$ws.Cells.Item(1564,5).Value2 = 'Normal'
$ws.Cells.Item(1564,6).Value2 = '6 a 8'
$ws.Cells.Item(1564,7).Value2 = 'Adolorido de una zona'
$ws.Cells.Item(1564,8).Value2 = '16 espalda'

$ws.Range('A1545:H1545').Copy($ws.Range('A1565:H1565'))
$ws.Cells.Item(1565,1).Value2 = 45880.40407255787
$ws.Cells.Item(1565,2).Value2 = 'Lunes'
$ws.Cells.Item(1565,3).Value2 = 'Ramón Juárez'
$ws.Cells.Item(1565,4).Value2 = 'Normal'
$ws.Cells.Item(1565,5).Value2 = 'Normal'
$ws.Cells.Item(1565,6).Value2 = '6 a 8'
$ws.Cells.Item(1565,7).Value2 = 'Normal'
$ws.Cells.Item(1565,8).Value2 = 'Nada'

$ws.Range('A1544:H1544').Copy($ws.Range('A1566:H1566'))
$ws.Cells.Item(1566,1).Value2 = 45880.41896456019
$ws.Cells.Item(1566,2).Value2 = 'Lunes'
$ws.Cells.Item(1566,3).Value2 = 'Miguel Vázquez'
$ws.Cells.Item(1566,4).Value2 = 'Normal'
$ws.Cells.Item(1566,5).Value2 = 'Normal'
$ws.Cells.Item(1566,6).Value2 = '6 a 8'
$ws.Cells.Item(1566,7).Value2 = 'Normal'
$ws.Cells.Item(1566,8).Value2 = 'Nada'

$ws.Range('A1545:H1545').Copy($ws.Range('A1567:H1567'))
$ws.Cells.Item(1567,1).Value2 = 45880.41930699074
$ws.Cells.Item(1567,2).Value2 = 'Lunes'
$ws.Cells.Item(1567,3).Value2 = 'Dagoberto Espinoza'
$ws.Cells.Item(1567,4).Value2 = 'Normal'
$ws.Cells.Item(1567,5).Value2 = 'Normal'
$ws.Cells.Item(1567,6).Value2 = '6 a 8'
$ws.Cells.Item(1567,7).Value2 = 'Normal'
$ws.Cells.Item(1567,8).Value2 = 'Nada'


# ---------------------------------------------------------------------
# 3) Row 1567 is now the new last row of the table, so it needs the
#    special darker bottom border that visually closes the table. Only
#    the bottom edge differs from a normal "odd" row -- the engine
#    already uses the lighter F8F9FA tone for top/left/right on plain
#    rows, so recoloring just the bottom edge to the dark purple
#    accent reproduces the "last row" look.
# ---------------------------------------------------------------------
$lastRow = $ws.Range('A1567:H1567')
$lastRow.Borders.Item(9).Color = 0x652F44   # xlEdgeBottom -> 442F65 (BGR)

# ---------------------------------------------------------------------
# 4) Grow the table / AutoFilter / hidden _FilterDatabase name so they
#    cover the newly added rows.
# ---------------------------------------------------------------------
$lo = $ws.ListObjects.Item(1)
while ($lo.Range.Rows.Count -lt 1567) {
    $lo.ListRows.Add() | Out-Null
}

foreach ($dn in $wb.Names) {
    if ($dn.Name -like '*_FilterDatabase*') {
        $dn.RefersTo = "='Respuestas de formulario 1'!`$A`$1:`$H`$1567"
    }
}

# ---------------------------------------------------------------------
# 5) Extend the two conditional-formatting rules (column D "Muy
#    cansado" and column G "adolorido") so their ranges keep pace with
#    the larger table, matching the sheet's existing +100 row buffer
#    beyond the table's actual extent.
# ---------------------------------------------------------------------
$rangeD = $ws.Range('D1:D1646')
$fcD = $rangeD.FormatConditions.Item(1)
$fcD.ModifyAppliesToRange($ws.Range('D1:D1667'))

$rangeG = $ws.Range('G2:G1646')
$fcG = $rangeG.FormatConditions.Item(1)
$fcG.ModifyAppliesToRange($ws.Range('G2:G1667'))

